$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title and link
$ws.Range("D9").Value = "문송합니다? 공송합니다, 컴송합니다 아냐?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/sorry-majors/#utm_source=rss&utm_medium=rss&utm_campaign=sorry-majors"

# Row 29: update title and link
$ws.Range("D29").Value = "[Flutter] Flutter에서 GraphQL을 잘 사용해보자!"
$ws.Range("E29").Value = "https://blog.promedius.ai/flutter-graphqleul-jal-sayonghaeboja/"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] Decision Transformer: Reinforcement Learning via Sequence Modeling"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1801&mod=document&pageid=1"
